$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" - refresh the handback status report with
# the latest run: two source files were regenerated with new GUIDs/hashes
# and new handoff/handback timestamps.
#   58095348-c0b4-43d8-8f04-d10ba197177c  ->  135675ad-6f70-4b39-950c-de98fae9371f
#   d7cb7150-4012-4530-9ce7-a12d547371e3  ->  ffffc6b7a6d2-5fef-4cf2-8d09-28028b1e74ce
# ---------------------------------------------------------------------------

$oldGuid1 = "58095348-c0b4-43d8-8f04-d10ba197177c"
$newGuid1 = "135675ad-6f70-4b39-950c-de98fae9371f"
$oldGuid2 = "d7cb7150-4012-4530-9ce7-a12d547371e3"
$newGuid2 = "ffffc6b7a6d2-5fef-4cf2-8d09-28028b1e74ce"

$oldHash1 = "37828abf4fdb20550175f1504d3114480c20be32"
$newHash  = "2bd1f3af0fbfa928fc0510428fc7bfd5a7a4c2f8"
$oldHash2 = "b9e076fbd1627c8bf8d4c2b66c007d60010ded50"

# -------------------- Overview sheet --------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = "e2e\$newGuid1.md"
$ws.Range("G2").Value = "2016-08-25 11:05:37"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "e2e\$newGuid2.md"
$ws.Range("G3").Value = "2016-08-25 11:05:37"

# -------------------- zh-cn sheet --------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-25 11:05:32"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-25 11:05:48"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-25 11:05:32"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-25 11:05:48"

# -------------------- de-de sheet --------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-25 11:05:37"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("K2").Value = "2016-08-25 11:05:55"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H3").Value = "2016-08-25 11:05:37"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("K3").Value = "2016-08-25 11:05:55"

# -------------------- Hyperlinks --------------------
# Rewrite each sheet's hyperlink display text to match the new file names
# while preserving the original target URLs (those did not change).

function Update-SheetHyperlinks($ws, $links) {
    $targets = @()
    foreach ($h in $ws.Hyperlinks) {
        $targets += $h.Address
    }
    $ws.Hyperlinks.Delete()
    for ($i = 0; $i -lt $links.Count; $i++) {
        $ws.Hyperlinks.Add($ws.Range($links[$i].Cell), $targets[$i], "", "", $links[$i].Display)
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
Update-SheetHyperlinks $wsOverview @(
    @{ Cell = "B2"; Display = "e2e\$newGuid1.md" },
    @{ Cell = "B3"; Display = "e2e\$newGuid2.md" }
)

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-SheetHyperlinks $wsZh @(
    @{ Cell = "A2"; Display = "$newGuid1.md" },
    @{ Cell = "I2"; Display = "$newGuid1.md" },
    @{ Cell = "A3"; Display = "$newGuid2.md" },
    @{ Cell = "I3"; Display = "$newGuid2.md" }
)

$wsDe = $wb.Worksheets.Item("de-de")
Update-SheetHyperlinks $wsDe @(
    @{ Cell = "A2"; Display = "$newGuid1.md" },
    @{ Cell = "I2"; Display = "$newGuid1.md" },
    @{ Cell = "A3"; Display = "$newGuid2.md" },
    @{ Cell = "I3"; Display = "$newGuid2.md" }
)
